$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.654227
$ws.Cells.Item(2, 8).Value = 4.962681
$ws.Cells.Item(2, 9).Value = 0.4107585939979205
$ws.Cells.Item(2, 10).Value = 0.4107585939979205
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 10.02989533333333
$ws.Cells.Item(2, 14).Value = 30.089686
$ws.Cells.Item(2, 15).Value = 0.8389720380165088
$ws.Cells.Item(2, 16).Value = 0.8389720380165089
$ws.Cells.Item(2, 17).Value = 16.591723667574
$ws.Cells.Item(2, 18).Value = 149.325513008166
$ws.Cells.Item(2, 19).Value = 0.3446149747392311
$ws.Cells.Item(2, 20).Value = 0.3446149747392311
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.654227
$ws.Cells.Item(3, 8).Value = 4.962681
$ws.Cells.Item(3, 9).Value = 0.4107585939979205
$ws.Cells.Item(3, 10).Value = 0.4107585939979205
$ws.Cells.Item(3, 15).Value = 0.08062742045192038
$ws.Cells.Item(3, 16).Value = 0.0806274204519204
$ws.Cells.Item(3, 17).Value = 1.594508302482
$ws.Cells.Item(3, 18).Value = 14.350574722338
$ws.Cells.Item(3, 19).Value = 0.03311840586251
$ws.Cells.Item(3, 20).Value = 0.03311840586251
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.654227
$ws.Cells.Item(4, 8).Value = 4.962681
$ws.Cells.Item(4, 9).Value = 0.4107585939979205
$ws.Cells.Item(4, 10).Value = 0.4107585939979205
$ws.Cells.Item(4, 13).Value = 0.911782
$ws.Cells.Item(4, 14).Value = 2.735346
$ws.Cells.Item(4, 15).Value = 0.07626795468388421
$ws.Cells.Item(4, 16).Value = 0.07626795468388423
$ws.Cells.Item(4, 17).Value = 1.508294402514
$ws.Cells.Item(4, 18).Value = 13.574649622626
$ws.Cells.Item(4, 19).Value = 0.0313277178330494
$ws.Cells.Item(4, 20).Value = 0.03132771783304941
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.654227
$ws.Cells.Item(5, 8).Value = 4.962681
$ws.Cells.Item(5, 9).Value = 0.4107585939979205
$ws.Cells.Item(5, 10).Value = 0.4107585939979205
$ws.Cells.Item(5, 13).Value = 0.049405
$ws.Cells.Item(5, 14).Value = 0.148215
$ws.Cells.Item(5, 15).Value = 0.004132586847686508
$ws.Cells.Item(5, 16).Value = 0.004132586847686508
$ws.Cells.Item(5, 17).Value = 0.08172708493499999
$ws.Cells.Item(5, 18).Value = 0.7355437644149999
$ws.Cells.Item(5, 19).Value = 0.001697495563130008
$ws.Cells.Item(5, 20).Value = 0.001697495563130009
$ws.Cells.Item(6, 9).Value = 0.3200015957958394
$ws.Cells.Item(6, 10).Value = 0.3200015957958394
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 10.02989533333333
$ws.Cells.Item(6, 14).Value = 30.089686
$ws.Cells.Item(6, 15).Value = 0.8389720380165088
$ws.Cells.Item(6, 16).Value = 0.8389720380165089
$ws.Cells.Item(6, 17).Value = 12.92578689334533
$ws.Cells.Item(6, 18).Value = 116.332082040108
$ws.Cells.Item(6, 19).Value = 0.2684723909933704
$ws.Cells.Item(6, 20).Value = 0.2684723909933705
$ws.Cells.Item(7, 9).Value = 0.3200015957958394
$ws.Cells.Item(7, 10).Value = 0.3200015957958394
$ws.Cells.Item(7, 15).Value = 0.08062742045192038
$ws.Cells.Item(7, 16).Value = 0.0806274204519204
$ws.Cells.Item(7, 19).Value = 0.02580090320951662
$ws.Cells.Item(7, 20).Value = 0.02580090320951663
$ws.Cells.Item(8, 9).Value = 0.3200015957958394
$ws.Cells.Item(8, 10).Value = 0.3200015957958394
$ws.Cells.Item(8, 13).Value = 0.911782
$ws.Cells.Item(8, 14).Value = 2.735346
$ws.Cells.Item(8, 15).Value = 0.07626795468388421
$ws.Cells.Item(8, 16).Value = 0.07626795468388423
$ws.Cells.Item(8, 17).Value = 1.175037169732
$ws.Cells.Item(8, 18).Value = 10.575334527588
$ws.Cells.Item(8, 19).Value = 0.02440586720692771
$ws.Cells.Item(8, 20).Value = 0.02440586720692772
$ws.Cells.Item(9, 9).Value = 0.3200015957958394
$ws.Cells.Item(9, 10).Value = 0.3200015957958394
$ws.Cells.Item(9, 13).Value = 0.049405
$ws.Cells.Item(9, 14).Value = 0.148215
$ws.Cells.Item(9, 15).Value = 0.004132586847686508
$ws.Cells.Item(9, 16).Value = 0.004132586847686508
$ws.Cells.Item(9, 17).Value = 0.06366950802999999
$ws.Cells.Item(9, 18).Value = 0.5730255722699999
$ws.Cells.Item(9, 19).Value = 0.00132243438602458
$ws.Cells.Item(9, 20).Value = 0.00132243438602458
$ws.Cells.Item(10, 7).Value = 0.8858993333333333
$ws.Cells.Item(10, 8).Value = 2.657698
$ws.Cells.Item(10, 9).Value = 0.2199763179924491
$ws.Cells.Item(10, 10).Value = 0.2199763179924491
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.02989533333333
$ws.Cells.Item(10, 14).Value = 30.089686
$ws.Cells.Item(10, 15).Value = 0.8389720380165088
$ws.Cells.Item(10, 16).Value = 0.8389720380165089
$ws.Cells.Item(10, 17).Value = 8.885477589203111
$ws.Cells.Item(10, 18).Value = 79.96929830282799
$ws.Cells.Item(10, 19).Value = 0.1845539798214927
$ws.Cells.Item(10, 20).Value = 0.1845539798214927
$ws.Cells.Item(11, 7).Value = 0.8858993333333333
$ws.Cells.Item(11, 8).Value = 2.657698
$ws.Cells.Item(11, 9).Value = 0.2199763179924491
$ws.Cells.Item(11, 10).Value = 0.2199763179924491
$ws.Cells.Item(11, 15).Value = 0.08062742045192038
$ws.Cells.Item(11, 16).Value = 0.0806274204519204
$ws.Cells.Item(11, 17).Value = 0.8539177768004443
$ws.Cells.Item(11, 18).Value = 7.685259991203999
$ws.Cells.Item(11, 19).Value = 0.01773612308024253
$ws.Cells.Item(11, 20).Value = 0.01773612308024254
$ws.Cells.Item(12, 7).Value = 0.8858993333333333
$ws.Cells.Item(12, 8).Value = 2.657698
$ws.Cells.Item(12, 9).Value = 0.2199763179924491
$ws.Cells.Item(12, 10).Value = 0.2199763179924491
$ws.Cells.Item(12, 13).Value = 0.911782
$ws.Cells.Item(12, 14).Value = 2.735346
$ws.Cells.Item(12, 15).Value = 0.07626795468388421
$ws.Cells.Item(12, 16).Value = 0.07626795468388423
$ws.Cells.Item(12, 17).Value = 0.8077470659453333
$ws.Cells.Item(12, 18).Value = 7.269723593507999
$ws.Cells.Item(12, 19).Value = 0.01677714385217581
$ws.Cells.Item(12, 20).Value = 0.01677714385217582
$ws.Cells.Item(13, 7).Value = 0.8858993333333333
$ws.Cells.Item(13, 8).Value = 2.657698
$ws.Cells.Item(13, 9).Value = 0.2199763179924491
$ws.Cells.Item(13, 10).Value = 0.2199763179924491
$ws.Cells.Item(13, 13).Value = 0.049405
$ws.Cells.Item(13, 14).Value = 0.148215
$ws.Cells.Item(13, 15).Value = 0.004132586847686508
$ws.Cells.Item(13, 16).Value = 0.004132586847686508
$ws.Cells.Item(13, 17).Value = 0.04376785656333333
$ws.Cells.Item(13, 18).Value = 0.3939107090699999
$ws.Cells.Item(13, 19).Value = 0.0009090712385381001
$ws.Cells.Item(13, 20).Value = 0.0009090712385381003
$ws.Cells.Item(14, 7).Value = 0.1983963333333333
$ws.Cells.Item(14, 8).Value = 0.595189
$ws.Cells.Item(14, 9).Value = 0.04926349221379096
$ws.Cells.Item(14, 10).Value = 0.04926349221379096
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 10.02989533333333
$ws.Cells.Item(14, 14).Value = 30.089686
$ws.Cells.Item(14, 15).Value = 0.8389720380165088
$ws.Cells.Item(14, 16).Value = 0.8389720380165089
$ws.Cells.Item(14, 17).Value = 1.989894457850444
$ws.Cells.Item(14, 18).Value = 17.909050120654
$ws.Cells.Item(14, 19).Value = 0.04133069246241461
$ws.Cells.Item(14, 20).Value = 0.04133069246241462
$ws.Cells.Item(15, 7).Value = 0.1983963333333333
$ws.Cells.Item(15, 8).Value = 0.595189
$ws.Cells.Item(15, 9).Value = 0.04926349221379096
$ws.Cells.Item(15, 10).Value = 0.04926349221379096
$ws.Cells.Item(15, 15).Value = 0.08062742045192038
$ws.Cells.Item(15, 16).Value = 0.0806274204519204
$ws.Cells.Item(15, 17).Value = 0.1912340934357778
$ws.Cells.Item(15, 18).Value = 1.721106840922
$ws.Cells.Item(15, 19).Value = 0.003971988299651229
$ws.Cells.Item(15, 20).Value = 0.00397198829965123
$ws.Cells.Item(16, 7).Value = 0.1983963333333333
$ws.Cells.Item(16, 8).Value = 0.595189
$ws.Cells.Item(16, 9).Value = 0.04926349221379096
$ws.Cells.Item(16, 10).Value = 0.04926349221379096
$ws.Cells.Item(16, 13).Value = 0.911782
$ws.Cells.Item(16, 14).Value = 2.735346
$ws.Cells.Item(16, 15).Value = 0.07626795468388421
$ws.Cells.Item(16, 16).Value = 0.07626795468388423
$ws.Cells.Item(16, 17).Value = 0.1808942055993333
$ws.Cells.Item(16, 18).Value = 1.628047850394
$ws.Cells.Item(16, 19).Value = 0.003757225791731292
$ws.Cells.Item(16, 20).Value = 0.003757225791731292
$ws.Cells.Item(17, 7).Value = 0.1983963333333333
$ws.Cells.Item(17, 8).Value = 0.595189
$ws.Cells.Item(17, 9).Value = 0.04926349221379096
$ws.Cells.Item(17, 10).Value = 0.04926349221379096
$ws.Cells.Item(17, 13).Value = 0.049405
$ws.Cells.Item(17, 14).Value = 0.148215
$ws.Cells.Item(17, 15).Value = 0.004132586847686508
$ws.Cells.Item(17, 16).Value = 0.004132586847686508
$ws.Cells.Item(17, 17).Value = 0.009801770848333332
$ws.Cells.Item(17, 18).Value = 0.08821593763499999
$ws.Cells.Item(17, 19).Value = 0.0002035856599938192
$ws.Cells.Item(17, 20).Value = 0.0002035856599938192
